$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Fill in the new time-registration rows (6-12) ---
$ws.Range("A6").Value = "Diverse projektrelaterede opgaver"
$ws.Range("C6").Value = 43963
$ws.Range("D6").Value = 0.35416666666666669
$ws.Range("E6").Value = 0.375
$ws.Range("F6").Value = 0.020833333333333332

$ws.Range("A7").Value = "Lavet mockups til UC01 og UC02"
$ws.Range("C7").Value = 43963
$ws.Range("D7").Value = 0.375
$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = 0.083333333333333329

$ws.Range("A8").Value = "Review af visionsdokument"
$ws.Range("C8").Value = 43963
$ws.Range("D8").Value = 0.54166666666666663
$ws.Range("E8").Value = 0.54861111111111105
$ws.Range("F8").Value = 0.020833333333333332

$ws.Range("A9").Value = "Lavet readme"
$ws.Range("C9").Value = 43963
$ws.Range("D9").Value = 0.54861111111111105
$ws.Range("E9").Value = 0.55555555555555558
$ws.Range("F9").Value = 0.0069444444444444441

$ws.Range("A10").Value = "Review af risikoanalyse og rettelser til metrikker"
$ws.Range("C10").Value = 43963
$ws.Range("D10").Value = 0.55902777777777779
$ws.Range("E10").Value = 0.625
$ws.Range("F10").Value = 0.041666666666666664

$ws.Range("A11").Value = "Rettet readme"
$ws.Range("C11").Value = 43963
$ws.Range("D11").Value = 0.625
$ws.Range("E11").Value = 0.63194444444444442
$ws.Range("F11").Value = 0.003472222222222222

$ws.Range("A12").Value = "Projektplan og gruppemøde"
$ws.Range("C12").Value = 43963
$ws.Range("D12").Value = 0.64583333333333337
$ws.Range("E12").Value = 0.6875
$ws.Range("F12").Value = 0.0625

# --- Apply the "estimate" column formatting (hh:mm, right+vcenter) for the new rows ---
$ws.Range("F6:F12").NumberFormat = "h:mm"
$ws.Range("F6:F12").HorizontalAlignment = -4152
$ws.Range("F6:F12").VerticalAlignment = -4108

# --- Re-align the remaining (still-empty) estimate cells ---
$ws.Range("F13:F24").HorizontalAlignment = -4152
$ws.Range("F25:F32").HorizontalAlignment = -4152

# --- Update the active selection to match the edited workbook ---
$ws.Range("A13").Select()

Write-Output "done"
